# Update gh-pages output data (generated at 456a3b4)
# Applies updated "想去人数" (interest count) values to 展览 / 演出 / 全部类型 sheets,
# and one "最低票价" correction on 展览!G12.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 86
$ws.Range("F5").Value  = 1675
$ws.Range("F6").Value  = 3268
$ws.Range("F7").Value  = 872
$ws.Range("F8").Value  = 2076
$ws.Range("F9").Value  = 1993
$ws.Range("F10").Value = 1029
$ws.Range("F11").Value = 359
$ws.Range("G12").Value = 128
$ws.Range("F16").Value = 15
$ws.Range("F17").Value = 74
$ws.Range("F18").Value = 87
$ws.Range("F19").Value = 1458
$ws.Range("F20").Value = 536
$ws.Range("F21").Value = 642
$ws.Range("F22").Value = 329
$ws.Range("F23").Value = 10898
$ws.Range("F24").Value = 11803
$ws.Range("F25").Value = 864
$ws.Range("F26").Value = 667
$ws.Range("F27").Value = 1854
$ws.Range("F28").Value = 155
$ws.Range("F29").Value = 460

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 35

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value  = 86
$ws.Range("F7").Value  = 1675
$ws.Range("F8").Value  = 3268
$ws.Range("F9").Value  = 872
$ws.Range("F10").Value = 2076
$ws.Range("F11").Value = 1993
$ws.Range("F12").Value = 1029
$ws.Range("F13").Value = 359
$ws.Range("F18").Value = 15
$ws.Range("F20").Value = 74
$ws.Range("F22").Value = 87
$ws.Range("F23").Value = 1458
$ws.Range("F24").Value = 536
$ws.Range("F25").Value = 642
$ws.Range("F26").Value = 329
$ws.Range("F27").Value = 10900
$ws.Range("F28").Value = 11803
$ws.Range("F29").Value = 864
$ws.Range("F30").Value = 667
$ws.Range("F31").Value = 1854
$ws.Range("F33").Value = 35
$ws.Range("F34").Value = 155
$ws.Range("F35").Value = 460
